$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell {
    param($Sheet, $Address, $Val)
    $c = $Sheet.Range($Address)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "43.006.07"
Set-TextCell $ws "E2" "  +0.41%  "
Set-TextCell $ws "D3" "2.301.13"
Set-TextCell $ws "E3" "  -0.13%  "
Set-TextCell $ws "E4" "  +0.23%  "
Set-TextCell $ws "D5" "311.28"
Set-TextCell $ws "E5" "  -2.65%  "
Set-TextCell $ws "D6" "104.73"
Set-TextCell $ws "E6" "  +0.08%  "
Set-TextCell $ws "E7" "  -1.15%  "
Set-TextCell $ws "E8" "  +0.26%  "
Set-TextCell $ws "D9" "0.604"
Set-TextCell $ws "E9" "  -0.94%  "
Set-TextCell $ws "D10" "39.95"
Set-TextCell $ws "E10" "  -0.57%  "
Set-TextCell $ws "D11" "0.0902"
Set-TextCell $ws "E11" "  -0.85%  "
Set-TextCell $ws "D12" "8.25"
Set-TextCell $ws "E12" "  -4.38%  "
Set-TextCell $ws "E13" "  +0.42%  "
Set-TextCell $ws "D14" "0.985"
Set-TextCell $ws "E14" "  +0.93%  "
Set-TextCell $ws "D15" "2.771.08"
Set-TextCell $ws "E15" "  +4.27%  "
Set-TextCell $ws "D16" "15.35"
Set-TextCell $ws "E16" "  -0.20%  "
Set-TextCell $ws "D17" "2.296.47"
Set-TextCell $ws "E17" "  -0.19%  "
Set-TextCell $ws "D18" "42.829.73"
Set-TextCell $ws "E18" "  +0.38%  "
Set-TextCell $ws "D19" "7.30"
Set-TextCell $ws "E19" "  -3.02%  "
Set-TextCell $ws "B20" "InternetComputer(DFINITY)"
Set-TextCell $ws "C20" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D20" "13.64"
Set-TextCell $ws "E20" "  +1.10%  "
Set-TextCell $ws "B21" "ShibaInu"
Set-TextCell $ws "C21" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D21" "0.0000104"
Set-TextCell $ws "E21" "  -1.58%  "
Set-TextCell $ws "D22" "73.38"
Set-TextCell $ws "E22" "  -0.64%  "
Set-TextCell $ws "D23" "3.46"
Set-TextCell $ws "E23" "  -3.19%  "
Set-TextCell $ws "D24" "269.11"
Set-TextCell $ws "E24" "  -1.02%  "
Set-TextCell $ws "D25" "2.18"
Set-TextCell $ws "E25" "  -3.86%  "
Set-TextCell $ws "E26" "  +0.47%  "
Set-TextCell $ws "D27" "10.85"
Set-TextCell $ws "E27" "  -0.82%  "
Set-TextCell $ws "D28" "7.21"
Set-TextCell $ws "E28" "  +15.80%  "
Set-TextCell $ws "D29" "2.29"
Set-TextCell $ws "E29" "  -1.18%  "
Set-TextCell $ws "D30" "22.33"
Set-TextCell $ws "E30" "  -1.49%  "
Set-TextCell $ws "D31" "36.10"
Set-TextCell $ws "E31" "  -5.22%  "
Set-TextCell $ws "D32" "164.42"
Set-TextCell $ws "E32" "  -0.93%  "
Set-TextCell $ws "D33" "0.0854"
Set-TextCell $ws "E33" "  -4.06%  "
Set-TextCell $ws "E34" "  +3.60%  "
Set-TextCell $ws "E35" "  -1.84%  "
Set-TextCell $ws "E36" "  -3.34%  "
Set-TextCell $ws "D37" "4.56"
Set-TextCell $ws "E37" "  -1.28%  "
Set-TextCell $ws "D38" "0.0348"
Set-TextCell $ws "E38" "  -2.20%  "
Set-TextCell $ws "E39" "  +3.12%  "
Set-TextCell $ws "E40" "  -2.47%  "
Set-TextCell $ws "D41" "110.24"
Set-TextCell $ws "E41" "  +10.81%  "
Set-TextCell $ws "E42" "  +0.99%  "
Set-TextCell $ws "D43" "70.83"
Set-TextCell $ws "E43" "  +0.00%  "
Set-TextCell $ws "E44" "  +0.19%  "
Set-TextCell $ws "E45" "  +0.50%  "
Set-TextCell $ws "D46" "12.31"
Set-TextCell $ws "E46" "  -0.23%  "
Set-TextCell $ws "D47" "1.736.14"
Set-TextCell $ws "E47" "  +8.78%  "
Set-TextCell $ws "D48" "110.70"
Set-TextCell $ws "E48" "  -3.61%  "
Set-TextCell $ws "D49" "77.60"
Set-TextCell $ws "E49" "  -5.79%  "
Set-TextCell $ws "B50" "THORChain"
Set-TextCell $ws "C50" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws "D50" "5.16"
Set-TextCell $ws "E50" "  -3.08%  "
Set-TextCell $ws "B51" "FraxShare"
Set-TextCell $ws "C51" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D51" "8.65"
Set-TextCell $ws "E51" "  -2.52%  "
